# West Bengal SoIB_summaries.xlsx update
# - "Trends Status", "Priority Status" and "Species qualification" sheets get
#   refreshed numbers (memory-optimized pipeline re-run).
# - The old "High Priority break-up" sheet is duplicated: the duplicate keeps
#   the old (pre-refresh) numbers and becomes "Major update - High Priority ",
#   while the original sheet is refreshed with new numbers and renamed to
#   "Interannual update - High Pri".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Trends Status" sheet (sheet1): refreshed Trend Status numbers
# ---------------------------------------------------------------------------
$trends = $wb.Worksheets.Item("Trends Status")

$trends.Range("B2").Value = 0
$trends.Range("C2").Value = 2
$trends.Range("D2").Value = 0
$trends.Range("E2").Value = 6.5

$trends.Range("B3").Value = 0
$trends.Range("C3").Value = 7
$trends.Range("D3").Value = 0
$trends.Range("E3").Value = 22.6

$trends.Range("B4").Value = 3
$trends.Range("C4").Value = 19
$trends.Range("D4").Value = 100
$trends.Range("E4").Value = 61.3

$trends.Range("E5").Value = 6.5

$trends.Range("B6").Value = 0
$trends.Range("C6").Value = 1
$trends.Range("D6").Value = 0
$trends.Range("E6").Value = 3.2

$trends.Range("B7").Value = 29
$trends.Range("C7").Value = 153

$trends.Range("B8").Value = 662
$trends.Range("C8").Value = 510

# ---------------------------------------------------------------------------
# 2. "Priority Status" sheet (sheet3): refreshed counts
# ---------------------------------------------------------------------------
$priority = $wb.Worksheets.Item("Priority Status")

$priority.Range("B2").Value = 103
$priority.Range("B3").Value = 286
$priority.Range("B4").Value = 554

# ---------------------------------------------------------------------------
# 3. "Species qualification" sheet (sheet4): refreshed counts + label tweak
# ---------------------------------------------------------------------------
$qualification = $wb.Worksheets.Item("Species qualification")

$qualification.Range("A2").Value = "SoIB Assessment"
$qualification.Range("B2").Value = 694

$qualification.Range("B3").Value = 32
$qualification.Range("C3").Value = 3

$qualification.Range("B4").Value = 184
$qualification.Range("C4").Value = 31

# ---------------------------------------------------------------------------
# 4. Duplicate "High Priority break-up" BEFORE touching its numbers, so the
#    duplicate preserves the old (pre-refresh) High Priority break-up figures.
#    Worksheet.Copy(Before, After) clones the sheet (values + formatting)
#    and is placed right after the source sheet, at the end of the workbook.
# ---------------------------------------------------------------------------
$breakup = $wb.Worksheets.Item("High Priority break-up")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$breakup.Copy($null, $lastSheet)

$majorUpdate = $wb.Worksheets.Item($wb.Worksheets.Count)
$majorUpdate.Name = "Major update - High Priority "

# ---------------------------------------------------------------------------
# 5. Refresh the original "High Priority break-up" sheet with the new
#    interannual-update numbers, then rename it.
# ---------------------------------------------------------------------------
$breakup.Range("B2").Value = 57
$breakup.Range("C2").Value = 55.3
$breakup.Range("D2").Value = 57
$breakup.Range("E2").Value = 72.2

$breakup.Range("B3").Value = 46
$breakup.Range("C3").Value = 44.7
$breakup.Range("D3").Value = 22
$breakup.Range("E3").Value = 27.8

$breakup.Name = "Interannual update - High Pri"

Write-Output "done"
